$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 100
$ws.Range("B4").Value = 8602.31615895128
$ws.Range("B5").Value = 87304.42145456493
$ws.Range("B6").Value = 1575
$ws.Range("B7").Value = 1421.172084527004
$ws.Range("B8").Value = 20217.60000000038
$ws.Range("B9").Value = 2930.691642071566
$ws.Range("B10").Value = 156795.5720740293
$ws.Range("B11").Value = 0.07876114942834694
$ws.Range("B12").Value = 0.2631509125128196
$ws.Range("B13").Value = 0.3499999999999949
$ws.Range("B14").Value = 0.9982380283505567
$ws.Range("B15").Value = 0.8954327058854333
